$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = "Done"
$ws.Range("G3").Value = "Done"
$ws.Range("I3").Value = "Done"
$ws.Range("J3").Value = "Done"
